# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the Leve profit sheets
# (columns H:N -- currentAveragePrice*, LevePrice*, LeveProfit*) per the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 125: Body over Mind
$ws.Range("H125").Value = 1398.75
$ws.Range("I125").Value = 1110.1666
$ws.Range("J125").Value = 2264.5
$ws.Range("K125").Value = 9991.499400000001
$ws.Range("L125").Value = 20380.5
$ws.Range("M125").Value = -7531.499400000001
$ws.Range("N125").Value = -25300.5

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 885.4286
$ws.Range("I137").Value = 885.4286
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2656.2858
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -106.2857999999997
$ws.Range("N137").Value = ""

# Row 138: All-night Crafting
$ws.Range("H138").Value = 5522.5
$ws.Range("J138").Value = 5640.625
$ws.Range("L138").Value = 16921.875
$ws.Range("N138").Value = -27201.875

$ws = $wb.Worksheets.Item("ARM")
# Row 28: 246 Kinds of Cheese
$ws.Range("H28").Value = 2066.5
$ws.Range("I28").Value = 2066.5
$ws.Range("K28").Value = 2066.5
$ws.Range("M28").Value = -1874.5

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1037.5
$ws.Range("I74").Value = 882.6
$ws.Range("J74").Value = 1812
$ws.Range("K74").Value = 882.6
$ws.Range("L74").Value = 1812
$ws.Range("M74").Value = -8.600000000000023
$ws.Range("N74").Value = -3560

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1037.5
$ws.Range("I77").Value = 882.6
$ws.Range("J77").Value = 1812
$ws.Range("K77").Value = 4413
$ws.Range("L77").Value = 9060
$ws.Range("M77").Value = -45
$ws.Range("N77").Value = -17796

# Row 99: Home Cooking
$ws.Range("H99").Value = 2066.5
$ws.Range("I99").Value = 2066.5
$ws.Range("K99").Value = 2066.5
$ws.Range("M99").Value = 928.5

# Row 125: The Incomplete Costume
$ws.Range("H125").Value = 90000.2
$ws.Range("J125").Value = 90000.2
$ws.Range("L125").Value = 90000.2
$ws.Range("N125").Value = -99840.2

$ws = $wb.Worksheets.Item("BSM")
# Row 44: You Spin Me Round
$ws.Range("H44").Value = 46332.332
$ws.Range("I44").Value = 24000
$ws.Range("J44").Value = 57498.5
$ws.Range("K44").Value = 24000
$ws.Range("L44").Value = 57498.5
$ws.Range("M44").Value = -23503
$ws.Range("N44").Value = -58492.5

# Row 106: Fire for Hire
$ws.Range("H106").Value = 25681.834
$ws.Range("J106").Value = 25681.834
$ws.Range("L106").Value = 25681.834
$ws.Range("N106").Value = -28205.834

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 182.8
$ws.Range("I7").Value = 163.25
$ws.Range("J7").Value = 200.84616
$ws.Range("K7").Value = 163.25
$ws.Range("L7").Value = 200.84616
$ws.Range("M7").Value = -50.25
$ws.Range("N7").Value = -426.84616

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 534.7778
$ws.Range("I22").Value = 591.7143
$ws.Range("J22").Value = 335.5
$ws.Range("K22").Value = 591.7143
$ws.Range("L22").Value = 335.5
$ws.Range("M22").Value = -241.7143
$ws.Range("N22").Value = -1035.5

# Row 35: Storm of Swords
$ws.Range("H35").Value = 4047.3333
$ws.Range("I35").Value = 4047.3333
$ws.Range("K35").Value = 4047.3333
$ws.Range("M35").Value = -3753.3333

# Row 45: A Tree Grew in Gridania
$ws.Range("H45").Value = 32000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 32000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 32000
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -33186

# Row 107: Built to Last
$ws.Range("H107").Value = 525.86664
$ws.Range("I107").Value = 391.3846
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 391.3846
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = 1528.6154
$ws.Range("N107").Value = -5240

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 565295.5600000001
$ws.Range("J141").Value = 565295.5600000001
$ws.Range("L141").Value = 565295.5600000001
$ws.Range("N141").Value = -575655.5600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 103.9
$ws.Range("I12").Value = 7
$ws.Range("J12").Value = 114.666664
$ws.Range("K12").Value = 21
$ws.Range("L12").Value = 343.999992
$ws.Range("M12").Value = 152
$ws.Range("N12").Value = -689.999992

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 3957.9385
$ws.Range("J55").Value = 4152.6333
$ws.Range("L55").Value = 12457.8999
$ws.Range("N55").Value = -12811.8999

# Row 69: Loving That Muffin Top
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""

# Row 72: Muffin of the Morn (L)
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""

# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 2195.75
$ws.Range("I80").Value = 1926.6666
$ws.Range("J80").Value = 3003
$ws.Range("K80").Value = 5779.9998
$ws.Range("L80").Value = 9009
$ws.Range("M80").Value = -4843.9998
$ws.Range("N80").Value = -10881

# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 2195.75
$ws.Range("I83").Value = 1926.6666
$ws.Range("J83").Value = 3003
$ws.Range("K83").Value = 17339.9994
$ws.Range("L83").Value = 27027
$ws.Range("M83").Value = -12659.9994
$ws.Range("N83").Value = -36387

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 173.94118
$ws.Range("J2").Value = 124.4
$ws.Range("L2").Value = 124.4
$ws.Range("N2").Value = -350.4

# Row 43: Get the Green Stuff
$ws.Range("H43").Value = 18836.285
$ws.Range("I43").Value = 6713.5
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 6713.5
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -6562.5
$ws.Range("N43").Value = -35302

# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 5010.25
$ws.Range("I46").Value = 5010.25
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 5010.25
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4854.25
$ws.Range("N46").Value = ""

# Row 57: Gold Is So Last Year
$ws.Range("H57").Value = 1036.6666
$ws.Range("I57").Value = 1036.6666
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 1036.6666
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -216.6666
$ws.Range("N57").Value = ""

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 8366.444
$ws.Range("I80").Value = 9928.286
$ws.Range("K80").Value = 9928.286
$ws.Range("M80").Value = -8930.286

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 8366.444
$ws.Range("I83").Value = 9928.286
$ws.Range("K83").Value = 49641.43
$ws.Range("M83").Value = -44649.43

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 2700
$ws.Range("I7").Value = 2700
$ws.Range("K7").Value = 2700
$ws.Range("M7").Value = -2588

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1205

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1393

# Row 126: Battered Books
$ws.Range("H126").Value = 2700
$ws.Range("I126").Value = 2700
$ws.Range("K126").Value = 8100
$ws.Range("M126").Value = -5630

$ws = $wb.Worksheets.Item("WVR")
# Row 21: Don't Trew So Hard
$ws.Range("H21").Value = 25000
$ws.Range("J21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("N21").Value = -25470

# Row 29: Getting Handsy
$ws.Range("H29").Value = 3000
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2710

# Row 35: Pantser Corps
$ws.Range("H35").Value = 25000
$ws.Range("J35").Value = 25000
$ws.Range("L35").Value = 25000
$ws.Range("N35").Value = -25580

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 3625.7778
$ws.Range("I126").Value = 1818.8572
$ws.Range("K126").Value = 5456.571599999999
$ws.Range("M126").Value = -2986.571599999999

Write-Host "Applied scheduled market-data refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."